# Update cryptos list — values refreshed as of the latest data pull.
# Price strings in column D are plain text (not numbers) in the source
# data, so for any cell whose new value would otherwise be auto-detected
# by Excel as a number (i.e. it has exactly one decimal separator), force
# the cell to Text format first so the literal string is preserved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.641.50"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.301.48"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.84"
$ws.Range("E5").Value = "  -2.68%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.48"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.97%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  -0.69%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.39"
$ws.Range("E10").Value = "  -1.73%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.56%  "

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.26"
$ws.Range("E12").Value = "  -3.40%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.16%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.988"
$ws.Range("E14").Value = "  +1.10%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.783.52"
$ws.Range("E15").Value = "  +4.48%  "

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.35"
$ws.Range("E16").Value = "  -0.36%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.301.18"
$ws.Range("E17").Value = "  -0.17%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.801.68"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -3.82%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -1.27%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.28"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23 - PancakeSwap
$ws.Range("E23").Value = "  -3.34%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.58"
$ws.Range("E24").Value = "  -0.68%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -1.78%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.81%  "

# Row 27 & 28 - Cosmos / Filecoin swap places, values updated
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.34"
$ws.Range("E27").Value = "  +18.08%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.88"
$ws.Range("E28").Value = "  -0.16%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -1.37%  "

# Row 30 - EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.27"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.23"
$ws.Range("E31").Value = "  -4.40%  "

# Row 32 - Monero
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.28"
$ws.Range("E32").Value = "  -0.23%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0857"
$ws.Range("E33").Value = "  -4.07%  "

# Row 34 & 35 - Stellar / WEMIXToken swap places, values updated
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").Value = "  +2.00%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  -1.28%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -3.65%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  -1.70%  "

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0347"
$ws.Range("E38").Value = "  -1.66%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +1.75%  "

# Row 40 - NEARProtocol
$ws.Range("E40").Value = "  -2.66%  "

# Row 41 - BitcoinSV
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "109.87"
$ws.Range("E41").Value = "  +11.39%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  -2.33%  "

# Row 43 - MultiversX
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.20"
$ws.Range("E43").Value = "  +1.38%  "

# Row 44 - Algorand
$ws.Range("E44").Value = "  +0.09%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.25%  "

# Row 46 - Celestia
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.25"
$ws.Range("E46").Value = "  -0.98%  "

# Row 47 - Maker
$ws.Range("D47").Value = "1.732.64"
$ws.Range("E47").Value = "  +6.66%  "

# Row 48 - Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.84"
$ws.Range("E48").Value = "  -4.58%  "

# Row 49 - ordi
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.35"
$ws.Range("E49").Value = "  -6.10%  "

# Row 50 - FraxShare
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.65"
$ws.Range("E50").Value = "  -2.79%  "

# Row 51 - THORChain
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.14"
$ws.Range("E51").Value = "  -2.98%  "
